$d = $word.ActiveDocument

# 1) Simple text replacements in the syllabus table (class topics column)
$d.Content.Find.Execute("Local classification and regression", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Irreducible and Reducible error", 2) | Out-Null

$d.Content.Find.Execute("Overfitting and resampling techniques", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Local methods", 2) | Out-Null

$d.Content.Find.Execute("Decision trees for classification and regression", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Model selection", 2) | Out-Null

$d.Content.Find.Execute("Ensemble methods for classification", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Decision trees for classification and regression; random forests?", 2) | Out-Null

# 2) Move "Assignment 1 due at the start of class" from row 5 (class 4) col 4
#    to row 6 (class 5) col 4, including its bold Arial 9pt formatting.
$t = $d.Tables.Item(2)

$srcCell = $t.Cell(5, 4)
$dstCell = $t.Cell(6, 4)

$srcRange = $srcCell.Range
$srcRange.MoveEnd(1, -1) | Out-Null   # exclude the trailing end-of-cell marker

$dstRange = $dstCell.Range
$dstRange.MoveEnd(1, -1) | Out-Null   # exclude the trailing end-of-cell marker

$dstRange.Text = "Assignment 1 due at the start of class"
$dstRange.Font.Name = "Arial"
$dstRange.Font.Bold = $true
$dstRange.Font.Size = 9

$srcRange.Text = ""
$srcRange.Font.Name = "Cambria"
$srcRange.Font.Bold = $false
